$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "A1BG"
$ws.Range("B2").Value = "ENSG00000121410.11"
$ws.Range("C2").Value = 6.45
$ws.Range("D2").Value = 18.61
$ws.Range("E2").Value = -1.396
$ws.Range("F2").Value = 0.00000000000000000000000000000401

$ws.Range("A3").Value = "AACSP1"
$ws.Range("B3").Value = "ENSG00000250420.8"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 2.535
$ws.Range("E3").Value = -1.822
$ws.Range("F3").Value = "4.1e-229"

$ws.Range("A4").Value = "AAGAB"
$ws.Range("B4").Value = "ENSG00000103591.12"
$ws.Range("C4").Value = 20.52
$ws.Range("D4").Value = 52.548
$ws.Range("E4").Value = -1.315
$ws.Range("F4").Value = 0.00000000000000000000000000000000000000000000000000000000000000398

$ws.Range("A5").Value = "AAMDC"
$ws.Range("B5").Value = "ENSG00000087884.14"
$ws.Range("C5").Value = 13.39
$ws.Range("D5").Value = 28.505
$ws.Range("E5").Value = -1.036
$ws.Range("F5").Value = 0.000000000000000000000000000000000000000000000000366

$ws.Range("A6").Value = "AARS"
$ws.Range("B6").Value = "ENSG00000090861.15"
$ws.Range("C6").Value = 31.441
$ws.Range("D6").Value = 205.329
$ws.Range("E6").Value = -2.669
$ws.Range("F6").Value = 0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000238

$ws.Range("A7").Value = "ABALON"
$ws.Range("B7").Value = "ENSG00000281376.1"
$ws.Range("C7").Value = 0.96
$ws.Range("D7").Value = 3.3
$ws.Range("E7").Value = -1.133
$ws.Range("F7").Value = 0.000000000000000000000000000000000000000000000000000000000000128

$ws.Range("A8").Value = "ABCA3"
$ws.Range("B8").Value = "ENSG00000167972.13"
$ws.Range("C8").Value = 2.38
$ws.Range("D8").Value = 7.385
$ws.Range("E8").Value = -1.311
$ws.Range("F8").Value = 0.00000000000000000000000264

$ws.Range("A9").Value = "ABCA8"
$ws.Range("B9").Value = "ENSG00000141338.13"
$ws.Range("C9").Value = 0.16
$ws.Range("D9").Value = 12.8
$ws.Range("E9").Value = -3.572
$ws.Range("F9").Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000231

$ws.Range("A10").Value = "ABCB10P1"
$ws.Range("B10").Value = "ENSG00000274099.1"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 2.75
$ws.Range("E10").Value = -1.907
$ws.Range("F10").Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000012

$ws.Range("A11").Value = "ABCB10P3"
$ws.Range("B11").Value = "ENSG00000261524.1"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 10.505
$ws.Range("E11").Value = -3.524
$ws.Range("F11").Value = "4.13e-251"

$ws.Range("A12").Value = "ABCB10P4"
$ws.Range("B12").Value = "ENSG00000260053.2"
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 25.325
$ws.Range("E12").Value = -4.718
$ws.Range("F12").Value = "1.42e-256"

$ws.Range("A13").Value = "ABCB6"
$ws.Range("B13").Value = "ENSG00000115657.12"
$ws.Range("C13").Value = 5.31
$ws.Range("D13").Value = 20.2
$ws.Range("E13").Value = -1.748
$ws.Range("F13").Value = 0.00000000000000000000000000000134

$ws.Range("A14").Value = "ABCB8"
$ws.Range("B14").Value = "ENSG00000197150.12"
$ws.Range("C14").Value = 18.02
$ws.Range("D14").Value = 45.984
$ws.Range("E14").Value = -1.305
$ws.Range("F14").Value = 0.00000000000000000000000000000000000000000000000000021

$ws.Range("A15").Value = "ABCC4"
$ws.Range("B15").Value = "ENSG00000125257.13"
$ws.Range("C15").Value = 9.07
$ws.Range("D15").Value = 20.554
$ws.Range("E15").Value = -1.098
$ws.Range("F15").Value = 0.000000000000000246

$ws.Range("A16").Value = "ABCF2"
$ws.Range("B16").Value = "ENSG00000033050.7"
$ws.Range("C16").Value = 24.609
$ws.Range("D16").Value = 54.57
$ws.Range("E16").Value = -1.118
$ws.Range("F16").Value = 0.000000000000000000000000000000000000000000000000000063

$ws.Range("A17").Value = "ABHD11"
$ws.Range("B17").Value = "ENSG00000106077.18"
$ws.Range("C17").Value = 7.67
$ws.Range("D17").Value = 21.41
$ws.Range("E17").Value = -1.37
$ws.Range("F17").Value = 0.00000000000000000000000000000000545

$ws.Range("A18").Value = "ABHD5"
$ws.Range("B18").Value = "ENSG00000011198.7"
$ws.Range("C18").Value = 13.85
$ws.Range("D18").Value = 38.019
$ws.Range("E18").Value = -1.394
$ws.Range("F18").Value = 0.000000000000000000000000000000000104

$ws.Range("A19").Value = "ABO"
$ws.Range("B19").Value = "ENSG00000175164.13"
$ws.Range("C19").Value = 1.52
$ws.Range("D19").Value = 10.93
$ws.Range("E19").Value = -2.243
$ws.Range("F19").Value = 0.0000000000103

$ws.Range("A20").Value = "ABT1"
$ws.Range("B20").Value = "ENSG00000146109.4"
$ws.Range("C20").Value = 9.66
$ws.Range("D20").Value = 24.48
$ws.Range("E20").Value = -1.257
$ws.Range("F20").Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000217

$ws.Range("A21").Value = "AC000041.8"
$ws.Range("B21").Value = "ENSG00000242156.1"
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 7.885
$ws.Range("E21").Value = -3.151
$ws.Range("F21").Value = 0.000000000000000000000000000000135
